# Updated cryptos list (GitHub Actions scrape refresh) - updates Price (D)
# and Volume(1h) (E) columns, plus three row reorderings (rows 12/13,
# 14/15) and a coin swap (row 51: MultiversX -> RocketPoolETH).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.269.91'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '2.216.22'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '301.55'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').Value = '88.76'
$ws.Range('E6').Value = '  -4.78%  '
$ws.Range('D7').Value = '0.552'
$ws.Range('E7').Value = '  -3.18%  '
$ws.Range('D9').Value = '0.492'
$ws.Range('E9').Value = '  -5.32%  '
$ws.Range('D10').Value = '33.38'
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').Value = '0.0776'
$ws.Range('E11').Value = '  -3.68%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.103'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '6.88'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.365.45'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.555.54'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '0.798'
$ws.Range('E16').Value = '  -3.30%  '
$ws.Range('D17').Value = '13.09'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').Value = '44.084.94'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').Value = '0.0₃0901'
$ws.Range('E19').Value = '  -6.31%  '
$ws.Range('D20').Value = '5.97'
$ws.Range('E20').Value = '  -5.54%  '
$ws.Range('D21').Value = '11.25'
$ws.Range('E21').Value = '  -6.03%  '
# D22/D26/D30/D34 are prices whose text ends in a trailing zero (e.g.
# "64.20"); a plain .Value assignment would be auto-parsed as a number
# and drop the trailing zero ("64.2"). Force text storage, then clear
# the number-format override so no extra style is left behind.
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.20'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').Value = '232.13'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('D24').Value = '2.86'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.74%  '
$ws.Range('E27').Value = '  +2.31%  '
$ws.Range('D28').Value = '9.37'
$ws.Range('E28').Value = '  -4.40%  '
$ws.Range('D29').Value = '36.01'
$ws.Range('E29').Value = '  -9.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.30'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').Value = '5.58'
$ws.Range('E31').Value = '  -4.65%  '
$ws.Range('D32').Value = '146.33'
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0750'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.25%  '
$ws.Range('D35').Value = '2.93'
$ws.Range('E35').Value = '  -3.72%  '
$ws.Range('E36').Value = '  -2.45%  '
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  -3.60%  '
$ws.Range('D38').Value = '1.73'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('D39').Value = '14.26'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').Value = '3.19'
$ws.Range('E40').Value = '  -7.59%  '
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('D42').Value = '0.0285'
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '1.735.54'
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('D45').Value = '1.68'
$ws.Range('E45').Value = '  +5.44%  '
$ws.Range('D46').Value = '78.16'
$ws.Range('E46').Value = '  -4.67%  '
$ws.Range('D47').Value = '0.179'
$ws.Range('E47').Value = '  -5.72%  '
$ws.Range('D48').Value = '94.44'
$ws.Range('E48').Value = '  -4.65%  '
$ws.Range('D49').Value = '4.65'
$ws.Range('E49').Value = '  -5.09%  '
$ws.Range('D50').Value = '66.21'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.437.59'
$ws.Range('E51').Value = '  -0.87%  '
